$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number + date range) ---
$ws.Range("A8").Characters(21,2).Text = "17"
$ws.Range("C9").Characters(27,9).Text = "4/22/2024"
$ws.Range("C9").Characters(47,9).Text = "4/28/2024"

# --- Data table updates (rows 15-28) ---
# Row 15
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E15").PasteSpecial(-4122)

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("F16").Value = 8
$ws.Range("H16").Value = 33.333333333333
$ws.Range("I16").Value = 26
$ws.Range("K16").Value = 8.333333333333
$ws.Range("L16").Value = 13.043478260869

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 20
$ws.Range("F17").Value = 15
$ws.Range("H17").Value = -25
$ws.Range("I17").Value = 68
$ws.Range("J17").Value = 78
$ws.Range("K17").Value = -12.820512820512
$ws.Range("L17").Value = -12.820512820512

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -14.285714285714
$ws.Range("I18").Value = 19
$ws.Range("J18").Value = 16
$ws.Range("K18").Value = 18.75
$ws.Range("L18").Value = -34.482758620689

# Row 19
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -40
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = -31.25
$ws.Range("I19").Value = 108
$ws.Range("J19").Value = 106
$ws.Range("K19").Value = 1.88679245283
$ws.Range("L19").Value = -30.322580645161

# Row 20
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -33.333333333333
$ws.Range("J20").Value = 29
$ws.Range("K20").Value = -48.275862068965

# Row 21
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = -10.526315789473
$ws.Range("F21").Value = 58
$ws.Range("G21").Value = 75
$ws.Range("H21").Value = -22.666666666666
$ws.Range("I21").Value = 239
$ws.Range("J21").Value = 257
$ws.Range("K21").Value = -7.003891050583
$ws.Range("L21").Value = -23.642172523961

# Row 23
$ws.Range("D23").Value = 1
$ws.Range("F24").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = -100
$ws.Range("E24").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("H23").Value = -100
$ws.Range("J23").Value = 8
$ws.Range("K23").Value = -75

# Row 24
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -39.393939393939
$ws.Range("F24").Value = 89
$ws.Range("G24").Value = 108
$ws.Range("H24").Value = -17.592592592592
$ws.Range("I24").Value = 490
$ws.Range("J24").Value = 390
$ws.Range("K24").Value = 25.641025641025
$ws.Range("L24").Value = 25

# Row 25
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 24
$ws.Range("E25").Value = -54.166666666666
$ws.Range("F25").Value = 55
$ws.Range("G25").Value = 73
$ws.Range("H25").Value = -24.657534246575
$ws.Range("I25").Value = 337
$ws.Range("J25").Value = 250
$ws.Range("K25").Value = 34.8
$ws.Range("L25").Value = 69.346733668341

# Row 26
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = 15
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 52
$ws.Range("G26").Value = 55
$ws.Range("H26").Value = -5.454545454545
$ws.Range("I26").Value = 207
$ws.Range("J26").Value = 180
$ws.Range("K26").Value = 15
$ws.Range("L26").Value = 30.188679245283

# Row 27
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("F24").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 16
$ws.Range("J28").Value = 23
$ws.Range("K28").Value = -30.434782608695
$ws.Range("L28").Value = 23.076923076923

$excel.CutCopyMode = 0
Write-Host "Edit script completed"